$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1) Remove the old "_GoBack" bookmark first (it used to sit at the end of
#    the totals paragraph, right after the final "10"). Doing this before
#    minting the new one avoids a transient duplicate-name collision.
# ---------------------------------------------------------------------------
$oldBookmark = $d.Bookmarks.Item("_GoBack")
$oldBookmark.Delete()

# ---------------------------------------------------------------------------
# 2) Title paragraph: "... Block - 869" -> "... Block - 50", and add the
#    "_GoBack" bookmark right after that run.
#
#    Both edits are applied with a single InsertXML call on the exact range
#    that held "869" so that (a) the sibling runs ("Block", " - ") - which
#    happen to carry identical rPr - are NOT coalesced into the rewritten
#    run, and (b) the bookmark tags land inside THIS paragraph (appending a
#    collapsed range's InsertXML at a paragraph-end position otherwise mints
#    a stray empty paragraph in this runtime).
# ---------------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$full = $titlePara.Range.Text
$idx = $full.IndexOf("869")
$start = $titlePara.Range.Start + $idx
$end = $start + 3
$numRange = $d.Range($start, $end)

# NB: the bookmark id below is deliberately "99" (not "0") - inserting a
# literal id="0" right after deleting the bookmark that used to own id="0"
# trips an internal id-collision path in this runtime that pollutes the
# package with unrelated namespace declarations. Using a throwaway id lets
# the engine renumber it (it comes out as id="0" in the saved file anyway,
# since it is the only bookmark left) without that side effect.
$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?><pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r w:rsidR="00E514AC"><w:rPr><w:rFonts w:ascii="Courier New" w:hAnsi="Courier New" w:cs="Courier New"/><w:b/><w:color w:val="1F3864" w:themeColor="accent5" w:themeShade="80"/><w:sz w:val="38"/><w:szCs w:val="38"/></w:rPr><w:t>50</w:t></w:r><w:bookmarkStart w:id="99" w:name="_GoBack"/><w:bookmarkEnd w:id="99"/></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$numRange.InsertXML($xml) | Out-Null
